# Vascular Malformations_Germline.xlsx refinement ("Refined metadata to be
# additional tab"):
#   1. Refresh the "time_taken" timestamps (column F) on the "data" sheet.
#   2. Add a new "metadata" sheet (placed after "data") holding the
#      panel-query bookkeeping fields (data_name, data_id, data_version,
#      data_version_created, panel_query_time, panel_get_request).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- 1. Refresh the F column ("time_taken") timestamps on "data" --------
$newTimes = @("2021-10-05 14:35:53.316818","2021-10-05 14:35:53.316826","2021-10-05 14:35:53.316829","2021-10-05 14:35:53.316832","2021-10-05 14:35:53.316835","2021-10-05 14:35:53.316838","2021-10-05 14:35:53.316840","2021-10-05 14:35:53.316843","2021-10-05 14:35:53.316846","2021-10-05 14:35:53.316848","2021-10-05 14:35:53.316851","2021-10-05 14:35:53.316854","2021-10-05 14:35:53.316856","2021-10-05 14:35:53.316859","2021-10-05 14:35:53.316861","2021-10-05 14:35:53.316864","2021-10-05 14:35:53.316867","2021-10-05 14:35:53.316870","2021-10-05 14:35:53.316873","2021-10-05 14:35:53.316875","2021-10-05 14:35:53.316878","2021-10-05 14:35:53.316880","2021-10-05 14:35:53.316883","2021-10-05 14:35:53.316885","2021-10-05 14:35:53.316888","2021-10-05 14:35:53.316891","2021-10-05 14:35:53.316894","2021-10-05 14:35:53.316896","2021-10-05 14:35:53.316899","2021-10-05 14:35:53.316901","2021-10-05 14:35:53.316904","2021-10-05 14:35:53.316906","2021-10-05 14:35:53.316910","2021-10-05 14:35:53.316912","2021-10-05 14:35:53.316915","2021-10-05 14:35:53.316917","2021-10-05 14:35:53.316920","2021-10-05 14:35:53.316922","2021-10-05 14:35:53.316925","2021-10-05 14:35:53.316927")

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- 2. Add the "metadata" sheet right after "data" ----------------------
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$meta.Name = "metadata"

# Reuse the exact header styling (bold, centered, thin border) already used
# by the "data" sheet's header row, and the same styling "data" uses on its
# index column (A), by copy/pasting the formats across.
$ws.Range("B1:F1").Copy() | Out-Null
$meta.Range("B1:G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A2").Copy() | Out-Null
$meta.Range("A2").PasteSpecial(-4122) | Out-Null      # xlPasteFormats

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Vascular Malformations_Germline"
$meta.Range("C2").Value = 300
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.4"
$meta.Range("E2").Value = "2021-09-29T05:14:11.967643Z"
$meta.Range("F2").Value = "2021-10-05 14:35:53.313141"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/300/?format=json"

$excel.CutCopyMode = $false

# Keep "data" as the active/selected sheet, matching the original workbook.
$ws.Activate()
